$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates: (cell address, new text value)
# All values are written as literal text (matching the original
# inline-string cell type) to avoid Excel auto-converting numeric-
# looking strings (e.g. "1.000", "235.27") into actual numbers.
$updates = @(
    @("D2", "30.646.49"),
    @("E2", "  +0.99%  "),
    @("D3", "1.868.96"),
    @("E3", "  +0.39%  "),
    @("D4", "1.000"),
    @("E4", "  +0.00%  "),
    @("D5", "235.27"),
    @("E5", "  +0.81%  "),
    @("E6", "  +0.01%  "),
    @("D7", "0.4698"),
    @("E7", "  -1.27%  "),
    @("D8", "0.2768"),
    @("E8", "  +0.77%  "),
    @("E9", "  -0.91%  "),
    @("D10", "18.04"),
    @("E10", "  +11.46%  "),
    @("D11", "1.864.94"),
    @("E11", "  +0.31%  "),
    @("D12", "0.07448"),
    @("E12", "  +0.22%  "),
    @("E13", "  +0.10%  "),
    @("D14", "85.38"),
    @("E14", "  -0.42%  "),
    @("D15", "0.6389"),
    @("E15", "  +1.21%  "),
    @("D16", "30.610.00"),
    @("E16", "  +1.01%  "),
    @("D17", "241.54"),
    @("E17", "  +2.87%  "),
    @("E18", "  +0.08%  "),
    @("D19", "12.87"),
    @("E19", "  +0.51%  "),
    @("D20", "0.000007405"),
    @("E20", "  +0.43%  "),
    @("D21", "1.0000"),
    @("E21", "  +0.09%  "),
    @("D22", "5.000"),
    @("E22", "  -1.58%  "),
    @("D23", "6.077"),
    @("E23", "  +1.38%  "),
    @("D24", "9.406"),
    @("E24", "  +1.38%  "),
    @("D25", "165.97"),
    @("E25", "  +0.04%  "),
    @("E26", "  +2.11%  "),
    @("D27", "1.894"),
    @("E27", "  +1.70%  "),
    @("E28", "  +1.94%  "),
    @("D29", "1.380"),
    @("E29", "  -0.04%  "),
    @("D30", "4.099"),
    @("E30", "  -2.63%  "),
    @("D31", "3.877"),
    @("E31", "  -1.21%  "),
    @("D32", "0.04937"),
    @("E32", "  +0.65%  "),
    @("D33", "1.155"),
    @("E33", "  +0.78%  "),
    @("D34", "0.7128"),
    @("E34", "  -1.07%  "),
    @("D35", "2.708"),
    @("E35", "  +0.42%  "),
    @("D36", "0.01912"),
    @("E36", "  +0.08%  "),
    @("D37", "2.699"),
    @("E37", "  +2.45%  "),
    @("D38", "0.8802"),
    @("E38", "  -2.48%  "),
    @("D39", "1.997"),
    @("E39", "  +0.82%  "),
    @("D40", "105.97"),
    @("E40", "  +0.23%  "),
    @("D41", "0.9998"),
    @("E41", "  -0.06%  "),
    @("D42", "0.4119"),
    @("E42", "  +0.50%  "),
    @("D43", "5.557"),
    @("E43", "  -0.07%  "),
    @("D44", "7.422"),
    @("E44", "  +5.15%  "),
    @("D45", "62.59"),
    @("E45", "  +2.58%  "),
    @("D46", "0.1236"),
    @("E46", "  +2.56%  "),
    @("B47", "EnergySwap"),
    @("C47", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"),
    @("D47", "8.704"),
    @("E47", "  -1.02%  "),
    @("B48", "Elrond"),
    @("C48", "https://coinranking.com/coin/omwkOTglq+elrond-egld"),
    @("D48", "33.70"),
    @("E48", "  +2.10%  "),
    @("D49", "0.05576"),
    @("E49", "  -0.28%  "),
    @("D50", "1.379"),
    @("E50", "  -1.58%  "),
    @("D51", "0.3719"),
    @("E51", "  +0.62%  ")
)

foreach ($pair in $updates) {
    $addr = $pair[0]
    $text = $pair[1]
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Write-Output "Applied $($updates.Count) cell updates"